$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04560231152239558
$ws.Range("H2").Value = 7.861300561477106
$ws.Range("I2").Value = 54.64886602885748
$ws.Range("G3").Value = 0.05480887362384868
$ws.Range("H3").Value = 12.54308134832489
$ws.Range("G4").Value = -0.009210837112246226
$ws.Range("H4").Value = -583.7889920768956
$ws.Range("G5").Value = 0.01540883072149391
$ws.Range("H5").Value = 305.8853102583673
$ws.Range("G6").Value = 0.05268921891400878
$ws.Range("H6").Value = 51.97536208828781
$ws.Range("G7").Value = -0.0001769237259479474
$ws.Range("H7").Value = -100.3326325068676
$ws.Range("G8").Value = -0.01203855901600541
$ws.Range("H8").Value = 36.03579554925334
$ws.Range("G9").Value = 0.00383017176652322
$ws.Range("H9").Value = 117.7916604704317
$ws.Range("G10").Value = -0.0779681464262218
$ws.Range("H10").Value = -7.244418418945194
$ws.Range("G11").Value = -0.05844105697944695
$ws.Range("H11").Value = 36.48473486704552
$ws.Range("G12").Value = -0.2448913669027963
$ws.Range("H12").Value = -0.1786518180766903
$ws.Range("G13").Value = -0.2662279137990329
$ws.Range("H13").Value = 3.124051709457758
$ws.Range("G14").Value = -0.06811065361976452
$ws.Range("H14").Value = -83.60055571595734
$ws.Range("G15").Value = -0.03830589927722763
$ws.Range("H15").Value = -10.16509001066514
$ws.Range("G16").Value = 0.127356010275682
$ws.Range("H16").Value = 1.634229069845331
$ws.Range("G17").Value = 0.1411099200003327
$ws.Range("H17").Value = 0.6127813251152966
$ws.Range("G18").Value = 0.1295946850868651
$ws.Range("H18").Value = 3.897903393182502
$ws.Range("G19").Value = 0.1285057946960539
$ws.Range("H19").Value = -3.5342034133129
$ws.Range("G20").Value = 0.05118338567372036
$ws.Range("H20").Value = 49.06575505159518
$ws.Range("G21").Value = 0.08507132971615736
$ws.Range("H21").Value = 46.57396810582471
$ws.Range("G22").Value = -0.08015290375337156
$ws.Range("H22").Value = -0.3818751891475089
$ws.Range("G23").Value = -0.08426563797050143
$ws.Range("H23").Value = -34.72492678251735
$ws.Range("G24").Value = 0.1236740116289681
$ws.Range("H24").Value = 4.706874362504264
$ws.Range("G25").Value = 0.1277987880720661
$ws.Range("H25").Value = 1.290631092178516
$ws.Range("G26").Value = 0.04434990064630216
$ws.Range("H26").Value = -10.77216300722759
$ws.Range("G27").Value = 0.07912491060091732
$ws.Range("H27").Value = -8.712457906771322
$ws.Range("G28").Value = -0.08525459836001245
$ws.Range("H28").Value = -34.06068095433106
$ws.Range("G29").Value = -0.07389216428936139
$ws.Range("H29").Value = -3.821635915781636
$ws.Range("G30").Value = 0.04734712397934603
$ws.Range("H30").Value = -25.68118643267752
$ws.Range("G31").Value = 0.03867978614674235
$ws.Range("H31").Value = -36.15138389623989
$ws.Range("G32").Value = 0.08952302705507752
$ws.Range("H32").Value = -8.895382441923781
$ws.Range("G33").Value = 0.09840413472291823
$ws.Range("H33").Value = 19.5898248198762
$ws.Range("G34").Value = 0.003465999333824717
$ws.Range("H34").Value = -86.69748376508096
$ws.Range("G35").Value = 0.004561685503989554
$ws.Range("H35").Value = 140.6926065528509
$ws.Range("G36").Value = 0.01497901514552191
$ws.Range("H36").Value = 2683.132304138926
$ws.Range("G37").Value = 0.01831168412873542
$ws.Range("H37").Value = 245.8602651316877
$ws.Range("G38").Value = 0.1247083146959879
$ws.Range("H38").Value = 16.27027024075836
$ws.Range("G39").Value = 0.1113234708348611
$ws.Range("H39").Value = 29.95546989835935
$ws.Range("G40").Value = 0.0222697908248077
$ws.Range("H40").Value = 649.759656088147
$ws.Range("G41").Value = 0.02450101693929442
$ws.Range("H41").Value = 63.37554981684399
$ws.Range("G42").Value = 0.1164157947907645
$ws.Range("H42").Value = 15.33724675549467
$ws.Range("G43").Value = 0.1043706321176854
$ws.Range("H43").Value = -13.12908262792253
$ws.Range("G44").Value = 0.02128126688405072
$ws.Range("H44").Value = -40.3678740187343
$ws.Range("G45").Value = 0.01482312736067157
$ws.Range("H45").Value = -9.450399166885612
$ws.Range("G46").Value = 0.06245818157168448
$ws.Range("H46").Value = 72.35301841330393
$ws.Range("G47").Value = 0.06025931220960556
$ws.Range("H47").Value = 19.46656868397554
$ws.Range("G48").Value = 0.0325803491511016
$ws.Range("H48").Value = -23.84441927489332
$ws.Range("G49").Value = 0.0571779234180126
$ws.Range("H49").Value = -17.7013061098494
$ws.Range("G50").Value = 0.02285432307763173
$ws.Range("H50").Value = 32.31466124896574
$ws.Range("G51").Value = 0.0336081113699283
$ws.Range("H51").Value = 72.61741219354971
$ws.Range("G52").Value = -0.1146752902331229
$ws.Range("H52").Value = -10.7754885511571
$ws.Range("G53").Value = -0.08409352872391211
$ws.Range("H53").Value = 8.945350861656296
$ws.Range("G54").Value = 0.07977332382905648
$ws.Range("H54").Value = 9.095214392778557
$ws.Range("G55").Value = 0.07394450847734348
$ws.Range("H55").Value = 19.35840243722351
$ws.Range("G56").Value = 0.03109048444348109
$ws.Range("H56").Value = -11.14266353056967
$ws.Range("G57").Value = 0.02866961826810242
$ws.Range("H57").Value = 396.5710616262593
$ws.Range("G58").Value = 0.04529116095698538
$ws.Range("H58").Value = 81.08781225687636
$ws.Range("G59").Value = 0.03455642447577283
$ws.Range("H59").Value = 45.93931588832823
$ws.Range("G60").Value = 0.02921688444278031
$ws.Range("H60").Value = -9.943000215885938
$ws.Range("G61").Value = 0.02249262547667658
$ws.Range("H61").Value = 77.69503565796779
$ws.Range("G62").Value = 0.05688905847556601
$ws.Range("H62").Value = -5.754487719326323
$ws.Range("G63").Value = 0.04651718184026379
$ws.Range("H63").Value = 42.73674447222295
$ws.Range("G64").Value = 0.03421533071850221
$ws.Range("H64").Value = -15.57227978131575
$ws.Range("G65").Value = 0.06042688437329032
$ws.Range("H65").Value = 7.785967970661295
$ws.Range("G66").Value = 0.1090907559321947
$ws.Range("H66").Value = 16.60712954432901
$ws.Range("G67").Value = 0.08925361875051105
$ws.Range("H67").Value = -22.68843333159854
$ws.Range("G68").Value = -0.02844111823196467
$ws.Range("H68").Value = 18.39103538736987
$ws.Range("G69").Value = -0.01177002466872856
$ws.Range("H69").Value = 44.53818550639994
$ws.Range("G70").Value = 0.08475393619452318
$ws.Range("H70").Value = -8.508894828605778
$ws.Range("G71").Value = 0.1040528231893123
$ws.Range("H71").Value = 14.08200460134278
$ws.Range("G72").Value = -0.0533195950210853
$ws.Range("H72").Value = 4.924002483775533
$ws.Range("G73").Value = -0.06612633715551289
$ws.Range("H73").Value = 10.35247370041895
$ws.Range("G74").Value = 0.1211774652428133
$ws.Range("H74").Value = 21.24195772581747
$ws.Range("G75").Value = 0.1045121705090663
$ws.Range("H75").Value = 7.29759164249395
$ws.Range("G76").Value = 0.005792477837151994
$ws.Range("H76").Value = -77.34711657925793
$ws.Range("G77").Value = 0.01300602464535132
$ws.Range("H77").Value = -7.817792644147237
$ws.Range("G78").Value = 0.100315789515285
$ws.Range("H78").Value = 56.06801576636925
$ws.Range("G79").Value = 0.09804389740097742
$ws.Range("H79").Value = 27.80461195061778
$ws.Range("G80").Value = -0.1495969118918831
$ws.Range("H80").Value = 9.668283864129066
$ws.Range("G81").Value = -0.1380216932368708
$ws.Range("H81").Value = 34.30416974342035
$ws.Range("G82").Value = 0.1257018018279704
$ws.Range("H82").Value = 9.596270396206551
$ws.Range("G83").Value = 0.2101594080440368
$ws.Range("H83").Value = 18.07888276255433
$ws.Range("G84").Value = 0.05984911152112442
$ws.Range("H84").Value = 151.0696898432162
$ws.Range("G85").Value = 0.07995278777763996
$ws.Range("H85").Value = 29.84420133853085
